$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(6047, 45863),
    @(5969, 45863.01041666666),
    @(5934, 45863.02083333334),
    @(5912, 45863.03125),
    @(5896, 45863.04166666666),
    @(5810, 45863.05208333334),
    @(5792, 45863.0625),
    @(5804, 45863.07291666666),
    @(5767, 45863.08333333334),
    @(5703, 45863.09375),
    @(5695, 45863.10416666666),
    @(5695, 45863.11458333334),
    @(5749, 45863.125),
    @(5753, 45863.13541666666),
    @(5734, 45863.14583333334),
    @(5687, 45863.15625),
    @(5748, 45863.16666666666),
    @(5807, 45863.17708333334),
    @(5767, 45863.1875),
    @(5799, 45863.19791666666),
    @(5875, 45863.20833333334),
    @(5903, 45863.21875),
    @(5981, 45863.22916666666),
    @(6084, 45863.23958333334),
    @(6247, 45863.25),
    @(6396, 45863.26041666666),
    @(6412, 45863.27083333334),
    @(6476, 45863.28125),
    @(6596, 45863.29166666666),
    @(6602, 45863.30208333334),
    @(6613, 45863.3125),
    @(6546, 45863.32291666666),
    @(6514, 45863.33333333334),
    @(6541, 45863.34375),
    @(6547, 45863.35416666666),
    @(6496, 45863.36458333334),
    @(6528, 45863.375),
    @(6555, 45863.38541666666)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
